# build min_em obj function
# Update the emission sheet with newly-computed solid_wood / sum_product /
# ecosystem / system values for periods 1 and 6-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 / F2 (period 1)
$ws.Range("E2").Value = -7.294328428059155
$ws.Range("F2").Value = -7.294328428059155

# E7 / F7 (period 6)
$ws.Range("E7").Value = -17.13070542992969
$ws.Range("F7").Value = -17.13070542992969

# Row 8 (period 7)
$ws.Range("B8").Value = 5.229794157038724
$ws.Range("D8").Value = 5.229794157038724
$ws.Range("E8").Value = 33.23912340812262
$ws.Range("F8").Value = 38.46891756516135

# Row 9 (period 8)
$ws.Range("B9").Value = 4.139652452119747
$ws.Range("D9").Value = 4.139652452119747
$ws.Range("E9").Value = -8.116932572908174
$ws.Range("F9").Value = -3.977280120788428

# Row 10 (period 9)
$ws.Range("B10").Value = 3.276748933086955
$ws.Range("D10").Value = 3.276748933086955
$ws.Range("E10").Value = -35.15785750188112
$ws.Range("F10").Value = -31.88110856879416

# Row 11 (period 10)
$ws.Range("B11").Value = 4.685633849059211
$ws.Range("D11").Value = 4.685633849059211
$ws.Range("E11").Value = -1.433666666666667
$ws.Range("F11").Value = 3.251967182392544
